$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(12, 8).Value = 312.25
$ws.Cells.Item(12, 9).Value = 249.66667
$ws.Cells.Item(12, 10).Value = 500
$ws.Cells.Item(12, 11).Value = 249.66667
$ws.Cells.Item(12, 12).Value = 500
$ws.Cells.Item(12, 13).Value = -79.66667000000001
$ws.Cells.Item(12, 14).Value = -840
$ws.Cells.Item(32, 8).Value = 2440
$ws.Cells.Item(32, 9).Value = 1400
$ws.Cells.Item(32, 10).Value = 4000
$ws.Cells.Item(32, 11).Value = 1400
$ws.Cells.Item(32, 12).Value = 4000
$ws.Cells.Item(32, 13).Value = -1074
$ws.Cells.Item(32, 14).Value = -4652
$ws.Cells.Item(61, 8).Value = 407.5
$ws.Cells.Item(61, 9).Value = 315
$ws.Cells.Item(61, 10).Value = 500
$ws.Cells.Item(61, 11).Value = 945
$ws.Cells.Item(61, 12).Value = 1500
$ws.Cells.Item(61, 13).Value = -773
$ws.Cells.Item(61, 14).Value = -1844
$ws.Cells.Item(86, 8).Value = 6563.9
$ws.Cells.Item(86, 9).Value = 4479.6665
$ws.Cells.Item(86, 10).Value = 7457.143
$ws.Cells.Item(86, 11).Value = 4479.6665
$ws.Cells.Item(86, 12).Value = 7457.143
$ws.Cells.Item(86, 13).Value = -3356.6665
$ws.Cells.Item(86, 14).Value = -9703.143
$ws.Cells.Item(89, 8).Value = 6563.9
$ws.Cells.Item(89, 9).Value = 4479.6665
$ws.Cells.Item(89, 10).Value = 7457.143
$ws.Cells.Item(89, 11).Value = 22398.3325
$ws.Cells.Item(89, 12).Value = 37285.715
$ws.Cells.Item(89, 13).Value = -16782.3325
$ws.Cells.Item(89, 14).Value = -48517.715
$ws.Cells.Item(92, 8).Value = 1331.091
$ws.Cells.Item(92, 9).Value = 1331.091
$ws.Cells.Item(92, 11).Value = 1331.091
$ws.Cells.Item(92, 13).Value = -83.09099999999989
$ws.Cells.Item(106, 8).Value = 2292.8333
$ws.Cells.Item(106, 9).Value = 2160.3125
$ws.Cells.Item(106, 10).Value = 3353
$ws.Cells.Item(106, 11).Value = 2160.3125
$ws.Cells.Item(106, 12).Value = 3353
$ws.Cells.Item(106, 13).Value = -1529.3125
$ws.Cells.Item(106, 14).Value = -4615
$ws.Cells.Item(138, 8).Value = 589457.5
$ws.Cells.Item(138, 9).Value = 1159.0952
$ws.Cells.Item(138, 10).Value = 822556.9
$ws.Cells.Item(138, 11).Value = 3477.2856
$ws.Cells.Item(138, 12).Value = 2467670.7
$ws.Cells.Item(138, 13).Value = 1662.7144
$ws.Cells.Item(138, 14).Value = -2477950.7
$ws.Cells.Item(141, 8).Value = 8277.267
$ws.Cells.Item(141, 9).Value = 8980.77
$ws.Cells.Item(141, 10).Value = 3704.5
$ws.Cells.Item(141, 11).Value = 26942.31
$ws.Cells.Item(141, 12).Value = 11113.5
$ws.Cells.Item(141, 13).Value = -21762.31
$ws.Cells.Item(141, 14).Value = -21473.5

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 2884.8164
$ws.Cells.Item(32, 9).Value = 2929.2917
$ws.Cells.Item(32, 11).Value = 2929.2917
$ws.Cells.Item(32, 13).Value = -2642.2917
$ws.Cells.Item(102, 8).Value = 20846890
$ws.Cells.Item(102, 9).Value = 33354154
$ws.Cells.Item(102, 10).Value = 1450
$ws.Cells.Item(102, 11).Value = 33354154
$ws.Cells.Item(102, 12).Value = 1450
$ws.Cells.Item(102, 13).Value = -33352532
$ws.Cells.Item(102, 14).Value = -4694
$ws.Cells.Item(110, 8).Value = 1188.7646
$ws.Cells.Item(110, 9).Value = 994.89655
$ws.Cells.Item(110, 11).Value = 994.89655
$ws.Cells.Item(110, 13).Value = 1050.10345
$ws.Cells.Item(132, 8).Value = 4992.727
$ws.Cells.Item(132, 9).Value = 5560.2856
$ws.Cells.Item(132, 10).Value = 3999.5
$ws.Cells.Item(132, 11).Value = 16680.8568
$ws.Cells.Item(132, 12).Value = 11998.5
$ws.Cells.Item(132, 13).Value = -14150.8568
$ws.Cells.Item(132, 14).Value = -17058.5
$ws.Cells.Item(133, 8).Value = 31217
$ws.Cells.Item(133, 10).Value = 31217
$ws.Cells.Item(133, 12).Value = 31217
$ws.Cells.Item(133, 14).Value = -36277

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(94, 8).Value = 41667868
$ws.Cells.Item(94, 9).Value = 83333800
$ws.Cells.Item(94, 10).Value = 1933
$ws.Cells.Item(94, 11).Value = 83333800
$ws.Cells.Item(94, 12).Value = 1933
$ws.Cells.Item(94, 13).Value = -83333349
$ws.Cells.Item(94, 14).Value = -2835
$ws.Cells.Item(99, 8).Value = 29413042
$ws.Cells.Item(99, 9).Value = 33334574
$ws.Cells.Item(99, 11).Value = 33334574
$ws.Cells.Item(99, 13).Value = -33333076
$ws.Cells.Item(107, 8).Value = 1495.0667
$ws.Cells.Item(107, 9).Value = 1201.3
$ws.Cells.Item(107, 11).Value = 1201.3
$ws.Cells.Item(107, 13).Value = 718.7
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 11208.637
$ws.Cells.Item(134, 9).Value = 2162
$ws.Cells.Item(134, 10).Value = 35333
$ws.Cells.Item(134, 11).Value = 6486
$ws.Cells.Item(134, 12).Value = 105999
$ws.Cells.Item(134, 13).Value = -3951
$ws.Cells.Item(134, 14).Value = -111069

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(7, 8).Value = 297.81818
$ws.Cells.Item(7, 9).Value = 303.25
$ws.Cells.Item(7, 11).Value = 303.25
$ws.Cells.Item(7, 13).Value = -190.25
$ws.Cells.Item(31, 8).Value = 1120.6786
$ws.Cells.Item(31, 9).Value = 884.1667
$ws.Cells.Item(31, 11).Value = 884.1667
$ws.Cells.Item(31, 13).Value = -589.1667
$ws.Cells.Item(34, 8).Value = 1120.6786
$ws.Cells.Item(34, 9).Value = 884.1667
$ws.Cells.Item(34, 11).Value = 884.1667
$ws.Cells.Item(34, 13).Value = -682.1667
$ws.Cells.Item(58, 8).Value = 832.8421
$ws.Cells.Item(58, 9).Value = 802.61536
$ws.Cells.Item(58, 11).Value = 802.61536
$ws.Cells.Item(58, 13).Value = -599.61536
$ws.Cells.Item(100, 8).Value = 99900
$ws.Cells.Item(100, 10).Value = 99900
$ws.Cells.Item(100, 12).Value = 99900
$ws.Cells.Item(100, 14).Value = -102064
$ws.Cells.Item(132, 8).Value = 6401.75
$ws.Cells.Item(132, 9).Value = 8609.286
$ws.Cells.Item(132, 10).Value = 3311.2
$ws.Cells.Item(132, 11).Value = 25827.858
$ws.Cells.Item(132, 12).Value = 9933.599999999999
$ws.Cells.Item(132, 13).Value = -23297.858
$ws.Cells.Item(132, 14).Value = -14993.6
$ws.Cells.Item(133, 8).Value = 63299.168
$ws.Cells.Item(133, 10).Value = 63299.168
$ws.Cells.Item(133, 12).Value = 63299.168
$ws.Cells.Item(133, 14).Value = -68359.16800000001
$ws.Cells.Item(134, 8).Value = 8131344
$ws.Cells.Item(134, 9).Value = 10102103
$ws.Cells.Item(134, 11).Value = 30306309
$ws.Cells.Item(134, 13).Value = -30303774
$ws.Cells.Item(136, 8).Value = 832.8421
$ws.Cells.Item(136, 9).Value = 802.61536
$ws.Cells.Item(136, 10).Value = 898.3333
$ws.Cells.Item(136, 11).Value = 2407.84608
$ws.Cells.Item(136, 12).Value = 2694.9999
$ws.Cells.Item(136, 13).Value = 142.1539199999997
$ws.Cells.Item(136, 14).Value = -7794.9999

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(68, 8).Value = 1474
$ws.Cells.Item(68, 9).Value = 749.1111
$ws.Cells.Item(68, 10).Value = 1784.6666
$ws.Cells.Item(68, 11).Value = 2247.3333
$ws.Cells.Item(68, 12).Value = 5353.9998
$ws.Cells.Item(68, 13).Value = -1436.3333
$ws.Cells.Item(68, 14).Value = -6975.9998
$ws.Cells.Item(71, 8).Value = 1474
$ws.Cells.Item(71, 9).Value = 749.1111
$ws.Cells.Item(71, 10).Value = 1784.6666
$ws.Cells.Item(71, 11).Value = 6741.9999
$ws.Cells.Item(71, 12).Value = 16061.9994
$ws.Cells.Item(71, 13).Value = -2685.9999
$ws.Cells.Item(71, 14).Value = -24173.9994
$ws.Cells.Item(86, 8).Value = 500
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 13).ClearContents()
$ws.Cells.Item(89, 8).Value = 500
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 13).ClearContents()
$ws.Cells.Item(107, 8).Value = 4615.16
$ws.Cells.Item(107, 9).Value = 605
$ws.Cells.Item(107, 10).Value = 7288.6
$ws.Cells.Item(107, 11).Value = 1815
$ws.Cells.Item(107, 12).Value = 21865.8
$ws.Cells.Item(107, 13).Value = 105
$ws.Cells.Item(107, 14).Value = -25705.8
$ws.Cells.Item(131, 8).Value = 23810950
$ws.Cells.Item(131, 9).Value = 200000720
$ws.Cells.Item(131, 10).Value = 1519.7838
$ws.Cells.Item(131, 11).Value = 600002160
$ws.Cells.Item(131, 12).Value = 4559.3514
$ws.Cells.Item(131, 13).Value = -599997120
$ws.Cells.Item(131, 14).Value = -14639.3514

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(70, 8).Value = 25005144
$ws.Cells.Item(70, 9).Value = 22731946
$ws.Cells.Item(70, 10).Value = 28577314
$ws.Cells.Item(70, 11).Value = 22731946
$ws.Cells.Item(70, 12).Value = 28577314
$ws.Cells.Item(70, 13).Value = -22731676
$ws.Cells.Item(70, 14).Value = -28577854
$ws.Cells.Item(73, 8).Value = 25005144
$ws.Cells.Item(73, 9).Value = 22731946
$ws.Cells.Item(73, 10).Value = 28577314
$ws.Cells.Item(73, 11).Value = 22731946
$ws.Cells.Item(73, 12).Value = 28577314
$ws.Cells.Item(73, 13).Value = -22731010
$ws.Cells.Item(73, 14).Value = -28579186
$ws.Cells.Item(126, 8).Value = 2083.9412
$ws.Cells.Item(126, 9).Value = 1647.909
$ws.Cells.Item(126, 10).Value = 2883.3333
$ws.Cells.Item(126, 11).Value = 4943.727000000001
$ws.Cells.Item(126, 12).Value = 8649.999899999999
$ws.Cells.Item(126, 13).Value = -2473.727000000001
$ws.Cells.Item(126, 14).Value = -13589.9999
$ws.Cells.Item(132, 8).Value = 3089.158
$ws.Cells.Item(132, 9).Value = 2933.111
$ws.Cells.Item(132, 10).Value = 3229.6
$ws.Cells.Item(132, 11).Value = 8799.332999999999
$ws.Cells.Item(132, 12).Value = 9688.799999999999
$ws.Cells.Item(132, 13).Value = -6269.332999999999
$ws.Cells.Item(132, 14).Value = -14748.8
$ws.Cells.Item(135, 8).Value = 49999
$ws.Cells.Item(135, 10).Value = 49999
$ws.Cells.Item(135, 12).Value = 49999
$ws.Cells.Item(135, 14).Value = -60139

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(41, 8).Value = 3518.5
$ws.Cells.Item(41, 9).Value = 2000
$ws.Cells.Item(41, 10).Value = 5037
$ws.Cells.Item(41, 11).Value = 2000
$ws.Cells.Item(41, 12).Value = 5037
$ws.Cells.Item(41, 13).Value = -1562
$ws.Cells.Item(41, 14).Value = -5913
$ws.Cells.Item(54, 8).Value = 10042
$ws.Cells.Item(54, 10).Value = 10042
$ws.Cells.Item(54, 12).Value = 10042
$ws.Cells.Item(54, 14).Value = -11330
$ws.Cells.Item(122, 8).Value = 31252308
$ws.Cells.Item(122, 9).Value = 35716496
$ws.Cells.Item(122, 11).Value = 107149488
$ws.Cells.Item(122, 13).Value = -107147038
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(58, 8).Value = 15000
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 10).Value = 15000
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 12).Value = 15000
$ws.Cells.Item(58, 13).ClearContents()
$ws.Cells.Item(58, 14).Value = -15616
